# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "G" column (header "K") previously held a "Strike#" style value; this
# script recalculates/rewrites it with the new "K" (strikeouts) values for
# every data row (rows 2-54), matching the regenerated save_data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..54 (row 52 recalculated to the same value it
# already had, so no visible change there).
$sVals = @(
    3, 2, 1, 0, 1, 2, 0, 1, 1, 0, 3, 1, 4, 0, 0, 1, 1, 2, 4, 2,
    2, 1, 0, 0, 3, 0, 1, 1, 1, 2, 3, 0, 1, 0, 0, 1, 1, 0, 1, 4,
    1, 1, 2, 1, 3, 1, 0, 1, 1, 1, 0, 0, 1
)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
